# Danh sach thanh vien - add "Tai khoan Github" column + new member row
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Insert a new column before column E ("So DT"), shifting E..H to F..I.
$ws.Range("E1").EntireColumn.Insert()
$ws.Range("E1").ColumnWidth = 26.14

# 2. Fill in the new column's values, in row-major (top-to-bottom) order so
#    the shared-string table is rebuilt in the same relative order as the
#    canonical file.
$ws.Range("E2").Value = "Tài khoản Github"

$ws.Range("E3").Value = "dtgianggithub"
$ws.Range("F3").Value = 1688452784

$ws.Range("E4").Value = "eooihic@gmail.com"

$ws.Range("E5").Value = "Thiện 1212381"

$ws.Range("E6").Value = "hongphuc4991@gmail.com"

$ws.Range("E7").Value = "nguyenhunghau.us@gmail.com"

# 3. New row for Nguyen Thanh Toan (MSSV 1212420)
$ws.Range("B8").Value = 1212420
$ws.Range("C8").Value = "Nguyễn Thành Toàn"
$ws.Range("D8").Value = "nguyenthanhtoan_94@yahoo.com"
$ws.Range("E8").Value = "nguyenthanhtoan"
$ws.Range("F8").Value = "0975178324"

# 4. Give new row-8 cells the normal sheet style (same as the rest of the data,
#    e.g. B3) so they're consistent with existing rows.
$ws.Range("B3").Copy()
$ws.Range("B8").PasteSpecial(-4122)
$ws.Range("B8").Value = 1212420

$ws.Range("D3").Copy()
$ws.Range("E8").PasteSpecial(-4122)
$ws.Range("E8").Value = "nguyenthanhtoan"

# 5. Apply the "pasted-from-web" grey highlight style (size 9, color #141823,
#    grey fill #F6F7F8, left aligned) to E5, then copy that exact resulting
#    style onto C8, D8 and F8 so they all reference the same new cellXf.
$e5 = $ws.Range("E5")
$e5.Font.Size = 9
$e5.Font.Color = 2299924
$e5.Interior.Color = 16316406
$e5.HorizontalAlignment = -4131

$e5.Copy()
$ws.Range("C8").PasteSpecial(-4122)
$ws.Range("D8").PasteSpecial(-4122)
$ws.Range("F8").PasteSpecial(-4122)

$ws.Range("C8").Value = "Nguyễn Thành Toàn"
$ws.Range("D8").Value = "nguyenthanhtoan_94@yahoo.com"
$ws.Range("F8").Value = "0975178324"

Write-Host "Edit applied"
